$wb = $excel.ActiveWorkbook

# --- Worksheet references (file order: test_suite, AddCustomerTest, OpenAccountTest) ---
$wsSuite = $wb.Worksheets.Item(1)   # test_suite          -> sheet1.xml
$wsAdd   = $wb.Worksheets.Item(2)   # AddCustomerTest      -> sheet2.xml
$wsOpen  = $wb.Worksheets.Item(3)   # OpenAccountTest      -> sheet3.xml

# ---------------------------------------------------------------------------
# 1) test_suite: flip the OpenAccountTest run-mode flag from N to Y
# ---------------------------------------------------------------------------
$wsSuite.Range("B3").Value = "Y"

# ---------------------------------------------------------------------------
# 2) AddCustomerTest: switch the data provider to a hashtable-style table
#    with two extra rows (Coto/Vladimir, MOS/CRACIUN) and a new runMode
#    column, plus re-cased headers.
#
#    NOTE: write the two brand-new data rows *before* rewriting the header
#    row so new shared strings are appended in the same order the source
#    workbook ended up with (Coto, Vladimir, MOS, CRACIUN, then the headers).
# ---------------------------------------------------------------------------
$wsAdd.Range("A4").Value = "Coto"
$wsAdd.Range("B4").Value = "Vladimir"
$wsAdd.Range("C4").Value = 123123
$wsAdd.Range("D4").Value = "Customer added successfully"
$wsAdd.Range("E4").Value = "N"

$wsAdd.Range("A5").Value = "MOS"
$wsAdd.Range("B5").Value = "CRACIUN"
$wsAdd.Range("C5").Value = 666
$wsAdd.Range("D5").Value = "Customer added successfully"
$wsAdd.Range("E5").Value = "Y"

$wsAdd.Range("A1").Value = "firstName"
$wsAdd.Range("B1").Value = "lastName"
$wsAdd.Range("C1").Value = "postCode"
$wsAdd.Range("D1").Value = "alertText"
$wsAdd.Range("E1").Value = "runMode"

$wsAdd.Range("E2").Value = "Y"
$wsAdd.Range("E3").Value = "Y"

# Approximate the bestFit-style column widths from the edited workbook.
# (ColumnWidth is internally quantized to 1/6-character steps here, so these
# inputs are chosen to land on the closest representable width to the
# target bestFit values of 9.1640625 / 13.1640625 / 8.5 / 25 / 8.5.)
$wsAdd.Columns.Item(1).ColumnWidth = 8.330729166666666
$wsAdd.Columns.Item(2).ColumnWidth = 12.330729166666666
$wsAdd.Columns.Item(3).ColumnWidth = 7.666666666666667
$wsAdd.Columns.Item(4).ColumnWidth = 24.166666666666668
$wsAdd.Columns.Item(5).ColumnWidth = 7.666666666666667

# ---------------------------------------------------------------------------
# 3) OpenAccountTest: re-case the alertText header to match the other sheet
# ---------------------------------------------------------------------------
$wsOpen.Range("C1").Value = "alertText"

# ---------------------------------------------------------------------------
# 4) Make OpenAccountTest the active tab (was test_suite) and move the
#    selection to C1.
# ---------------------------------------------------------------------------
$wsOpen.Activate()
$wsOpen.Range("C1").Select()
